$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Page_1")
$ws2 = $wb.Worksheets.Item("Page_2")

$oldText1 = '"Part Number - Can be found on the top right position of the page"'
$newText1 = '"Part Number - Can be found on the center right position of the page"'

$oldText2 = '"Duxford Range Part Number Description Dimensions Power Lumens Colour Temp. - Can be found on the bottom right position of the page"'
$newText2 = '"Multi-Wattage Tri-Colour and Single Colour 4000K Retrofit Gear Trays - Can be found on the middle right position of the page"'

for ($r = 2; $r -le 13; $r++) {
    $cell1 = $ws1.Range("G$r")
    if ($cell1.Value2 -eq $oldText1) {
        $cell1.Value = $newText1
    }

    $cell2 = $ws2.Range("G$r")
    if ($cell2.Value2 -eq $oldText2) {
        $cell2.Value = $newText2
    }
}
